# Edit: change the text "Data" to "Data Employee" inside the "What?" bullet
# of the "Who? / What? / Power BI" summary text box on slide 2.
#
# The original run (a:r/a:t) contains exactly the word "Data"; we locate it
# as a whole-word match inside the shape's TextRange and replace just that
# run's characters, leaving every other run / formatting run untouched.

$p = $ppt.ActivePresentation

$target = "Data"
$replacement = "Data Employee"
$found = $false
$alnum = "abcdefghijklmnopqrstuvwxyzABCDEFGHIJKLMNOPQRSTUVWXYZ0123456789"

for ($si = 1; $si -le $p.Slides.Count -and -not $found; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count -and -not $found; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }
        if (-not $shape.TextFrame.HasText) { continue }

        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        $searchFrom = 0

        while ($true) {
            $idx = $full.IndexOf($target, $searchFrom)
            if ($idx -lt 0) { break }

            $beforeOk = $true
            if ($idx -gt 0) {
                $prevChar = $full.Substring($idx - 1, 1)
                $beforeOk = -not $alnum.Contains($prevChar)
            }

            $afterPos = $idx + $target.Length
            $afterOk = $true
            if ($afterPos -lt $full.Length) {
                $nextChar = $full.Substring($afterPos, 1)
                $afterOk = -not $alnum.Contains($nextChar)
            }

            if ($beforeOk -and $afterOk) {
                $sub = $tr.Characters($idx + 1, $target.Length)
                if ($sub.Text -eq $target) {
                    $sub.Text = $replacement
                    $found = $true
                }
                break
            }

            $searchFrom = $idx + 1
        }
    }
}

if (-not $found) {
    Write-Output "WARNING: target text 'Data' not found"
} else {
    Write-Output "OK: replaced 'Data' with 'Data Employee'"
}
